# ============================================================================
# edit.ps1 - Chocobo_Profits.xlsx scheduled-runner update
#
# Refreshes cached market-board derived columns (currentAveragePrice*,
# LevePrice*/LeveProfit*) across all 8 job sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR) of the Chocobo_Profits leve-crafting workbook with freshly
# pulled price data. Columns H-N per sheet:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
# ============================================================================

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 563.44446  # H2: was 545.9
$ws.Cells.Item(2, 9).Value = 508.875  # I2: was 484.33334
$ws.Cells.Item(2, 10).Value = 1000  # J2: was 1100
$ws.Cells.Item(2, 11).Value = 508.875  # K2: was 484.33334
$ws.Cells.Item(2, 12).Value = 1000  # L2: was 1100
$ws.Cells.Item(2, 13).Value = -395.875  # M2: was -371.33334
$ws.Cells.Item(2, 14).Value = -1226  # N2: was -1326
$ws.Cells.Item(94, 8).Value = 1535.1666  # H94: was 1410.5714
$ws.Cells.Item(94, 9).Value = 1535.1666  # I94: was 1410.5714
$ws.Cells.Item(94, 11).Value = 1535.1666  # K94: was 1410.5714
$ws.Cells.Item(94, 13).Value = -1084.1666  # M94: was -959.5714
$ws.Cells.Item(116, 8).Value = 509415.44  # H116: was 565711.7
$ws.Cells.Item(116, 9).Value = 1432357.6  # I116: was 2004200.8
$ws.Cells.Item(116, 11).Value = 1432357.6  # K116: was 2004200.8
$ws.Cells.Item(116, 13).Value = -1428915.6  # M116: was -2000758.8
$ws.Cells.Item(129, 8).Value = 933.2632  # H129: was 922.2
$ws.Cells.Item(129, 10).Value = 937.8723  # J129: was 926.46466
$ws.Cells.Item(129, 12).Value = 2813.6169  # L129: was 2779.39398
$ws.Cells.Item(129, 14).Value = -12813.6169  # N129: was -12779.39398
$ws.Cells.Item(137, 8).Value = 2792.639  # H137: was 3149.5
$ws.Cells.Item(137, 9).Value = 1509.8182  # I137: was 1572.1904
$ws.Cells.Item(137, 10).Value = 4808.5  # J137: was 6829.8887
$ws.Cells.Item(137, 11).Value = 4529.4546  # K137: was 4716.5712
$ws.Cells.Item(137, 12).Value = 14425.5  # L137: was 20489.6661
$ws.Cells.Item(137, 13).Value = -1979.4546  # M137: was -2166.5712
$ws.Cells.Item(137, 14).Value = -19525.5  # N137: was -25589.6661
$ws.Cells.Item(138, 8).Value = 2418.26  # H138: was 2204.82
$ws.Cells.Item(138, 9).Value = 909.04346  # I138: was 752.5333000000001
$ws.Cells.Item(138, 10).Value = 2869.065  # J138: was 2827.2285
$ws.Cells.Item(138, 11).Value = 2727.13038  # K138: was 2257.5999
$ws.Cells.Item(138, 12).Value = 8607.195  # L138: was 8481.6855
$ws.Cells.Item(138, 13).Value = 2412.86962  # M138: was 2882.4001
$ws.Cells.Item(138, 14).Value = -18887.195  # N138: was -18761.6855

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5010.029  # H32: was 6043.3037
$ws.Cells.Item(32, 9).Value = 4254.8394  # I32: was 5233.711
$ws.Cells.Item(32, 10).Value = 8263.154  # J32: was 9355.272000000001
$ws.Cells.Item(32, 11).Value = 4254.8394  # K32: was 5233.711
$ws.Cells.Item(32, 12).Value = 8263.154  # L32: was 9355.272000000001
$ws.Cells.Item(32, 13).Value = -3967.8394  # M32: was -4946.711
$ws.Cells.Item(32, 14).Value = -8837.154  # N32: was -9929.272000000001
$ws.Cells.Item(61, 8).Value = 1700.9  # H61: was 1642.3
$ws.Cells.Item(61, 9).Value = 1101.6666  # I61: was 1002.875
$ws.Cells.Item(61, 10).Value = 2599.75  # J61: was 4200
$ws.Cells.Item(61, 11).Value = 1101.6666  # K61: was 1002.875
$ws.Cells.Item(61, 12).Value = 2599.75  # L61: was 4200
$ws.Cells.Item(61, 13).Value = -889.6666  # M61: was -790.875
$ws.Cells.Item(61, 14).Value = -3023.75  # N61: was -4624
$ws.Cells.Item(74, 8).Value = 4069.2593  # H74: was 4727.087
$ws.Cells.Item(74, 9).Value = 3886.0833  # I74: was 4605.95
$ws.Cells.Item(74, 11).Value = 3886.0833  # K74: was 4605.95
$ws.Cells.Item(74, 13).Value = -3012.0833  # M74: was -3731.95
$ws.Cells.Item(77, 8).Value = 4069.2593  # H77: was 4727.087
$ws.Cells.Item(77, 9).Value = 3886.0833  # I77: was 4605.95
$ws.Cells.Item(77, 11).Value = 19430.4165  # K77: was 23029.75
$ws.Cells.Item(77, 13).Value = -15062.4165  # M77: was -18661.75
$ws.Cells.Item(133, 8).Value = 30795  # H133: was 30930
$ws.Cells.Item(133, 10).Value = 30795  # J133: was 30930
$ws.Cells.Item(133, 12).Value = 30795  # L133: was 30930
$ws.Cells.Item(133, 14).Value = -35855  # N133: was -35990
$ws.Cells.Item(136, 8).Value = 1700.9  # H136: was 1642.3
$ws.Cells.Item(136, 9).Value = 1101.6666  # I136: was 1002.875
$ws.Cells.Item(136, 10).Value = 2599.75  # J136: was 4200
$ws.Cells.Item(136, 11).Value = 3304.9998  # K136: was 3008.625
$ws.Cells.Item(136, 12).Value = 7799.25  # L136: was 12600
$ws.Cells.Item(136, 13).Value = -754.9998000000001  # M136: was -458.625
$ws.Cells.Item(136, 14).Value = -12899.25  # N136: was -17700

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5037.625  # H20: was 7919.96
$ws.Cells.Item(20, 9).Value = 1416.862  # I20: was 2526.8
$ws.Cells.Item(20, 10).Value = 14583.272  # J20: was 16009.7
$ws.Cells.Item(20, 11).Value = 1416.862  # K20: was 2526.8
$ws.Cells.Item(20, 12).Value = 14583.272  # L20: was 16009.7
$ws.Cells.Item(20, 13).Value = -1169.862  # M20: was -2279.8
$ws.Cells.Item(20, 14).Value = -15077.272  # N20: was -16503.7
$ws.Cells.Item(86, 8).Value = 2344.111  # H86: was 1669.5555
$ws.Cells.Item(86, 9).Value = 1500  # I86: was 1232
$ws.Cells.Item(86, 10).Value = 2449.625  # J86: was 2019.6
$ws.Cells.Item(86, 11).Value = 1500  # K86: was 1232
$ws.Cells.Item(86, 12).Value = 2449.625  # L86: was 2019.6
$ws.Cells.Item(86, 13).Value = -377  # M86: was -109
$ws.Cells.Item(86, 14).Value = -4695.625  # N86: was -4265.6
$ws.Cells.Item(89, 8).Value = 2344.111  # H89: was 1669.5555
$ws.Cells.Item(89, 9).Value = 1500  # I89: was 1232
$ws.Cells.Item(89, 10).Value = 2449.625  # J89: was 2019.6
$ws.Cells.Item(89, 11).Value = 7500  # K89: was 6160
$ws.Cells.Item(89, 12).Value = 12248.125  # L89: was 10098
$ws.Cells.Item(89, 13).Value = -1884  # M89: was -544
$ws.Cells.Item(89, 14).Value = -23480.125  # N89: was -21330
$ws.Cells.Item(134, 8).Value = 1970.7715  # H134: was 1928.3611
$ws.Cells.Item(134, 10).Value = 2708.9333  # J134: was 2567.375
$ws.Cells.Item(134, 12).Value = 8126.7999  # L134: was 7702.125
$ws.Cells.Item(134, 14).Value = -13196.7999  # N134: was -12772.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8476399  # H31: was 10418709
$ws.Cells.Item(31, 9).Value = 949.0857  # I31: was 964.625
$ws.Cells.Item(31, 10).Value = 20836430  # J31: was 31254196
$ws.Cells.Item(31, 11).Value = 949.0857  # K31: was 964.625
$ws.Cells.Item(31, 12).Value = 20836430  # L31: was 31254196
$ws.Cells.Item(31, 13).Value = -654.0857  # M31: was -669.625
$ws.Cells.Item(31, 14).Value = -20837020  # N31: was -31254786
$ws.Cells.Item(34, 8).Value = 8476399  # H34: was 10418709
$ws.Cells.Item(34, 9).Value = 949.0857  # I34: was 964.625
$ws.Cells.Item(34, 10).Value = 20836430  # J34: was 31254196
$ws.Cells.Item(34, 11).Value = 949.0857  # K34: was 964.625
$ws.Cells.Item(34, 12).Value = 20836430  # L34: was 31254196
$ws.Cells.Item(34, 13).Value = -747.0857  # M34: was -762.625
$ws.Cells.Item(34, 14).Value = -20836834  # N34: was -31254600
$ws.Cells.Item(57, 8).Value = 41125.332  # H57: was 49999.8
$ws.Cells.Item(57, 9).Value = 10000  # I57: was 0
$ws.Cells.Item(57, 10).Value = 45016  # J57: was 49999.8
$ws.Cells.Item(57, 11).Value = 10000  # K57: was 0
$ws.Cells.Item(57, 12).Value = 45016  # L57: was 49999.8
$ws.Cells.Item(57, 13).Value = -9440  # M57: newly added
$ws.Cells.Item(57, 14).Value = -46136  # N57: was -51119.8
$ws.Cells.Item(59, 8).Value = 50000  # H59: was 37626
$ws.Cells.Item(59, 10).Value = 50000  # J59: was 37626
$ws.Cells.Item(59, 12).Value = 50000  # L59: was 37626
$ws.Cells.Item(59, 14).Value = -52290  # N59: was -39916
$ws.Cells.Item(62, 8).Value = 7166.6665  # H62: was 7500
$ws.Cells.Item(62, 10).Value = 5750  # J62: was 5000
$ws.Cells.Item(62, 12).Value = 5750  # L62: was 5000
$ws.Cells.Item(62, 14).Value = -6998  # N62: was -6248
$ws.Cells.Item(65, 8).Value = 7166.6665  # H65: was 7500
$ws.Cells.Item(65, 10).Value = 5750  # J65: was 5000
$ws.Cells.Item(65, 12).Value = 28750  # L65: was 25000
$ws.Cells.Item(65, 14).Value = -34990  # N65: was -31240
$ws.Cells.Item(134, 8).Value = 4041.7026  # H134: was 4842.8
$ws.Cells.Item(134, 9).Value = 4538.654  # I134: was 5490.048
$ws.Cells.Item(134, 10).Value = 2867.0908  # J134: was 3332.5557
$ws.Cells.Item(134, 11).Value = 13615.962  # K134: was 16470.144
$ws.Cells.Item(134, 12).Value = 8601.2724  # L134: was 9997.667099999999
$ws.Cells.Item(134, 13).Value = -11080.962  # M134: was -13935.144
$ws.Cells.Item(134, 14).Value = -13671.2724  # N134: was -15067.6671
$ws.Cells.Item(138, 8).Value = 41945  # H138: was 42733
$ws.Cells.Item(138, 10).Value = 41945  # J138: was 42733
$ws.Cells.Item(138, 12).Value = 41945  # L138: was 42733
$ws.Cells.Item(138, 14).Value = -52225  # N138: was -53013
$ws.Cells.Item(140, 8).Value = 76782  # H140: was 71846.47
$ws.Cells.Item(140, 10).Value = 76782  # J140: was 71846.47
$ws.Cells.Item(140, 12).Value = 76782  # L140: was 71846.47
$ws.Cells.Item(140, 14).Value = -87142  # N140: was -82206.47
$ws.Cells.Item(141, 8).Value = 32175  # H141: was 33037.5
$ws.Cells.Item(141, 10).Value = 32175  # J141: was 33037.5
$ws.Cells.Item(141, 12).Value = 32175  # L141: was 33037.5
$ws.Cells.Item(141, 14).Value = -42535  # N141: was -43397.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 13158830  # H131: was 15152453
$ws.Cells.Item(131, 10).Value = 1033.9395  # J131: was 1054.8928
$ws.Cells.Item(131, 12).Value = 3101.8185  # L131: was 3164.6784
$ws.Cells.Item(131, 14).Value = -13181.8185  # N131: was -13244.6784

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 7500  # H6: was 15000
$ws.Cells.Item(6, 10).Value = 7500  # J6: was 15000
$ws.Cells.Item(6, 12).Value = 7500  # L6: was 15000
$ws.Cells.Item(6, 14).Value = -7726  # N6: was -15226
$ws.Cells.Item(14, 8).Value = 148428600  # H14: was 171500130
$ws.Cells.Item(14, 9).Value = 148428600  # I14: was 171500130
$ws.Cells.Item(14, 11).Value = 148428600  # K14: was 171500130
$ws.Cells.Item(14, 13).Value = -148428432  # M14: was -171499962
$ws.Cells.Item(16, 8).Value = 7500  # H16: was 15000
$ws.Cells.Item(16, 10).Value = 7500  # J16: was 15000
$ws.Cells.Item(16, 12).Value = 7500  # L16: was 15000
$ws.Cells.Item(16, 14).Value = -8000  # N16: was -15500
$ws.Cells.Item(20, 8).Value = 0  # H20: was 25000
$ws.Cells.Item(20, 10).Value = 0  # J20: was 25000
$ws.Cells.Item(20, 12).Value = 0  # L20: was 25000
$ws.Cells.Item(20, 14).ClearContents()  # N20: was -25490, now blank
$ws.Cells.Item(21, 8).Value = 11166.333  # H21: was 11500
$ws.Cells.Item(21, 10).Value = 11749.5  # J21: was 12250
$ws.Cells.Item(21, 12).Value = 11749.5  # L21: was 12250
$ws.Cells.Item(21, 14).Value = -12095.5  # N21: was -12596
$ws.Cells.Item(30, 8).Value = 11166.333  # H30: was 11500
$ws.Cells.Item(30, 10).Value = 11749.5  # J30: was 12250
$ws.Cells.Item(30, 12).Value = 11749.5  # L30: was 12250
$ws.Cells.Item(30, 14).Value = -11959.5  # N30: was -12460
$ws.Cells.Item(80, 8).Value = 83335380  # H80: was 62502280
$ws.Cells.Item(80, 10).Value = 0  # J80: was 3000
$ws.Cells.Item(80, 12).Value = 0  # L80: was 3000
$ws.Cells.Item(80, 14).ClearContents()  # N80: was -4996, now blank
$ws.Cells.Item(83, 8).Value = 83335380  # H83: was 62502280
$ws.Cells.Item(83, 10).Value = 0  # J83: was 3000
$ws.Cells.Item(83, 12).Value = 0  # L83: was 15000
$ws.Cells.Item(83, 14).ClearContents()  # N83: was -24984, now blank

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 667.4211  # H16: was 1275.9
$ws.Cells.Item(16, 9).Value = 719.2308  # I16: was 1325.2858
$ws.Cells.Item(16, 10).Value = 555.1667  # J16: was 1160.6666
$ws.Cells.Item(16, 11).Value = 719.2308  # K16: was 1325.2858
$ws.Cells.Item(16, 12).Value = 555.1667  # L16: was 1160.6666
$ws.Cells.Item(16, 13).Value = -549.2308  # M16: was -1155.2858
$ws.Cells.Item(16, 14).Value = -895.1667  # N16: was -1500.6666
$ws.Cells.Item(68, 8).Value = 720.45715  # H68: was 901.8378
$ws.Cells.Item(68, 9).Value = 653.41174  # I68: was 716.7059
$ws.Cells.Item(68, 11).Value = 653.41174  # K68: was 716.7059
$ws.Cells.Item(68, 13).Value = 95.58825999999999  # M68: was 32.29409999999996
$ws.Cells.Item(71, 8).Value = 720.45715  # H71: was 901.8378
$ws.Cells.Item(71, 9).Value = 653.41174  # I71: was 716.7059
$ws.Cells.Item(71, 11).Value = 3267.0587  # K71: was 3583.5295
$ws.Cells.Item(71, 13).Value = 476.9413  # M71: was 160.4704999999999
$ws.Cells.Item(132, 8).Value = 8255.666999999999  # H132: was 2806.4138
$ws.Cells.Item(132, 9).Value = 3501  # I132: was 1203.9546
$ws.Cells.Item(132, 10).Value = 9614.143  # J132: was 7842.7144
$ws.Cells.Item(132, 11).Value = 10503  # K132: was 3611.8638
$ws.Cells.Item(132, 12).Value = 28842.429  # L132: was 23528.1432
$ws.Cells.Item(132, 13).Value = -7973  # M132: was -1081.8638
$ws.Cells.Item(132, 14).Value = -33902.429  # N132: was -28588.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 337.61905  # H113: was 356.10526
$ws.Cells.Item(113, 9).Value = 270.66666  # I113: was 284.25
$ws.Cells.Item(113, 10).Value = 387.83334  # J113: was 408.36365
$ws.Cells.Item(113, 11).Value = 811.9999799999999  # K113: was 852.75
$ws.Cells.Item(113, 12).Value = 1163.50002  # L113: was 1225.09095
$ws.Cells.Item(113, 13).Value = 1358.00002  # M113: was 1317.25
$ws.Cells.Item(113, 14).Value = -5503.500019999999  # N113: was -5565.09095
$ws.Cells.Item(132, 8).Value = 5953608.5  # H132: was 6290674.5
$ws.Cells.Item(132, 9).Value = 803.8163500000001  # I132: was 853.6
$ws.Cells.Item(132, 10).Value = 47623240  # J132: was 41670916
$ws.Cells.Item(132, 11).Value = 2411.44905  # K132: was 2560.8
$ws.Cells.Item(132, 12).Value = 142869720  # L132: was 125012748
$ws.Cells.Item(132, 13).Value = 118.5509499999998  # M132: was -30.80000000000018
$ws.Cells.Item(132, 14).Value = -142874780  # N132: was -125017808
$ws.Cells.Item(136, 8).Value = 2501.5789  # H136: was 2340.634
$ws.Cells.Item(136, 9).Value = 782.4  # I136: was 730.9286
$ws.Cells.Item(136, 11).Value = 2347.2  # K136: was 2192.7858
$ws.Cells.Item(136, 13).Value = 202.8000000000002  # M136: was 357.2142000000003
